$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 8.521337333333333
$ws.Range("N2").Value = 25.564012
$ws.Range("O2").Value = 0.2943426187002489
$ws.Range("P2").Value = 0.2943426187002489
$ws.Range("Q2").Value = 3.368797096902222
$ws.Range("R2").Value = 30.31917387212
$ws.Range("S2").Value = 0.162900937995511
$ws.Range("T2").Value = 0.162900937995511

# Row 3
$ws.Range("O3").Value = 0.1683364841626613
$ws.Range("P3").Value = 0.1683364841626613
$ws.Range("S3").Value = 0.09316412040517331
$ws.Range("T3").Value = 0.09316412040517334

# Row 4
$ws.Range("O4").Value = 0.5373208971370899
$ws.Range("P4").Value = 0.53732089713709
$ws.Range("S4").Value = 0.2973748026525505
$ws.Range("T4").Value = 0.2973748026525506

# Row 5
$ws.Range("M5").Value = 8.521337333333333
$ws.Range("N5").Value = 25.564012
$ws.Range("O5").Value = 0.2943426187002489
$ws.Range("P5").Value = 0.2943426187002489
$ws.Range("Q5").Value = 2.718218555514222
$ws.Range("R5").Value = 24.463966999628
$ws.Range("S5").Value = 0.1314416807047379
$ws.Range("T5").Value = 0.1314416807047379

# Row 6
$ws.Range("O6").Value = 0.1683364841626613
$ws.Range("P6").Value = 0.1683364841626613
$ws.Range("S6").Value = 0.07517236375748797
$ws.Range("T6").Value = 0.07517236375748798

# Row 7
$ws.Range("O7").Value = 0.5373208971370899
$ws.Range("P7").Value = 0.53732089713709
$ws.Range("S7").Value = 0.2399460944845394
$ws.Range("T7").Value = 0.2399460944845394
